# edit.ps1
# 1) Switch the three tables (slides 14, 15, 16) from the custom
#    "Table_0" style {9E8A5834-06BA-4847-A3DF-6E91A4004D17} to the
#    built-in "No Style, Table Grid" style
#    {29D4F1F3-7678-49BF-B9C9-F40F7F7D3F62}.
# 2) Re-colour the deck's primary theme (theme1.xml, used by the slide
#    master / all slides) from the "Integral / Red Violet" palette to
#    the stock "Office" palette (the palette that used to live in the
#    notes-master theme, theme2.xml).

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------
$newTableStyleId = "{29D4F1F3-7678-49BF-B9C9-F40F7F7D3F62}"

foreach ($slideIdx in @(14, 15, 16)) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newTableStyleId, $true)
    }
}

# --- 2) Theme colours --------------------------------------------------
function HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order for Slide.ThemeColorScheme.Colors(): Dark1, Light1, Dark2,
# Light2, Accent1-6, Hyperlink, FollowedHyperlink.
$officeThemeHex = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = HexToRgbInt($officeThemeHex[$i - 1])
}
